$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell H1: "Other parameter values:" (appended as new shared string)
# Copy the header style (bold/border/centered) from the existing G1 header cell
# so H1 matches the other header cells (B1:G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Other parameter values:"

# Size column H like the other data columns (target stored width = 24)
$ws.Range("H1").ColumnWidth = 23.2

# Clear the previous H2:H3 selection, leaving the default top-left cell selected
$ws.Range("A1").Select()
